$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 402 - this shifts the existing rows 402..504
# down to 403..505, matching the dimension growing from A1:R504 to A1:R505.
$ws.Rows.Item(402).Insert()

# Populate the newly inserted row 402 with the new weekly record.
$ws.Cells.Item(402, 1).Value = 4
$ws.Cells.Item(402, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(402, 3).Value = "Los Lagos"
$ws.Cells.Item(402, 4).Value = 45135
$ws.Cells.Item(402, 5).Value = 10
$ws.Cells.Item(402, 6).Value = 100112017
$ws.Cells.Item(402, 7).Value = "Apio"
$ws.Cells.Item(402, 8).Value = "Americana (o)"
$ws.Cells.Item(402, 9).Value = "Primera"
$ws.Cells.Item(402, 10).Value = 30
$ws.Cells.Item(402, 11).Value = 11000
$ws.Cells.Item(402, 12).Value = 11000
$ws.Cells.Item(402, 13).Value = 11000
$ws.Cells.Item(402, 14).Value = "`$/docena de matas"
$ws.Cells.Item(402, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(402, 16).Value = 1833
$ws.Cells.Item(402, 17).Value = 6
$ws.Cells.Item(402, 18).Value = "Hortaliza"
